$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9329097690173285
$ws.Range("C2").Value = 0.05228061303763809
$ws.Range("D2").Value = 0.1156676683936873
$ws.Range("E2").Value = 0.06005454773046726
$ws.Range("F2").Value = 2.075524740652597
$ws.Range("I2").Value = 1.688420095166613
$ws.Range("K2").Value = 0.7479508682363019
$ws.Range("L2").Value = 0.2249561431623661
$ws.Range("M2").Value = 0.2387367315177471
$ws.Range("N2").Value = 3.034829201778081
$ws.Range("B3").Value = 0.8966713425648436
$ws.Range("C3").Value = 0.04550233385448621
$ws.Range("D3").Value = 0.1157624700084838
$ws.Range("E3").Value = 0.06020361776105077
$ws.Range("F3").Value = 2.06421932893889
$ws.Range("I3").Value = 1.688000342013609
$ws.Range("K3").Value = 0.7065628561795165
$ws.Range("L3").Value = 0.2222859337750123
$ws.Range("M3").Value = 0.2320865127834075
$ws.Range("N3").Value = 3.050060697382982
$ws.Range("B4").Value = 0.8749476440333126
$ws.Range("C4").Value = 0.04133518365787836
$ws.Range("D4").Value = 0.1158200975665409
$ws.Range("E4").Value = 0.06030932923139698
$ws.Range("F4").Value = 2.058230358761762
$ws.Range("I4").Value = 1.688397463660664
$ws.Range("K4").Value = 0.6815773114576018
$ws.Range("L4").Value = 0.2207555003845272
$ws.Range("M4").Value = 0.2281344809301373
$ws.Range("N4").Value = 3.060176053919903
$ws.Range("B5").Value = 0.86622775164426
$ws.Range("C5").Value = 0.03963563088287003
$ws.Range("D5").Value = 0.1158434391262482
$ws.Range("E5").Value = 0.06035598249798912
$ws.Range("F5").Value = 2.056029287262149
$ws.Range("I5").Value = 1.688723973268502
$ws.Range("K5").Value = 0.6715028810066883
$ws.Range("L5").Value = 0.2201592984920069
$ws.Range("M5").Value = 0.2265570411517146
$ws.Range("N5").Value = 3.064489997482568
$ws.Range("B6").Value = 0.8647878407410303
$ws.Range("C6").Value = 0.03935333415734021
$ws.Range("D6").Value = 0.1158473065081065
$ws.Range("E6").Value = 0.06036394541826429
$ws.Range("F6").Value = 2.055678263633055
$ws.Range("I6").Value = 1.688788136566572
$ws.Range("K6").Value = 0.6698365165329108
$ws.Range("L6").Value = 0.220061959209751
$ws.Range("M6").Value = 0.2262971057268679
$ws.Range("N6").Value = 3.06521791065785
$ws.Range("B7").Value = 0.8748295069701726
$ws.Range("C7").Value = 0.04131226869255045
$ws.Range("D7").Value = 0.1158204129297395
$ws.Range("E7").Value = 0.06030994392685685
$ws.Range("F7").Value = 2.058199704766324
$ws.Range("I7").Value = 1.688401200293171
$ws.Range("K7").Value = 0.6814410093191441
$ws.Range("L7").Value = 0.2207473485614173
$ws.Range("M7").Value = 0.2281130731641774
$ws.Range("N7").Value = 3.060233456538789
$ws.Range("B8").Value = 0.9203055263999147
$ws.Range("C8").Value = 0.04994450224847924
$ws.Range("D8").Value = 0.1157004789457989
$ws.Range("E8").Value = 0.06010300888564046
$ws.Range("F8").Value = 2.071428879898548
$ws.Range("I8").Value = 1.688139452461634
$ws.Range("K8").Value = 0.7335916972905068
$ws.Range("L8").Value = 0.22401283162948
$ws.Range("M8").Value = 0.2364165245222765
$ws.Range("N8").Value = 3.039922598852954
$ws.Range("B9").Value = 1.013660706941948
$ws.Range("C9").Value = 0.06683506623942037
$ws.Range("D9").Value = 0.115460497481763
$ws.Range("E9").Value = 0.05980936469564391
$ws.Range("F9").Value = 2.104936121109688
$ws.Range("I9").Value = 1.6928234064912
$ws.Range("K9").Value = 0.839250430535202
$ws.Range("L9").Value = 0.2312812078828586
$ws.Range("M9").Value = 0.253739939928586
$ws.Range("N9").Value = 3.006150687704675
$ws.Range("B10").Value = 1.084799269779126
$ws.Range("C10").Value = 0.07922924020277833
$ws.Range("D10").Value = 0.1152810079933015
$ws.Range("E10").Value = 0.05966152085901477
$ws.Range("F10").Value = 2.13418016930774
$ws.Range("I10").Value = 1.699437477044697
$ws.Range("K10").Value = 0.9189607810573648
$ws.Range("L10").Value = 0.2371483412926949
$ws.Range("M10").Value = 0.2671022574942299
$ws.Range("N10").Value = 2.985034759134194
$ws.Range("B11").Value = 1.117717361209884
$ws.Range("C11").Value = 0.08486599471871159
$ws.Range("D11").Value = 0.1151986118552717
$ws.Range("E11").Value = 0.05960890761057946
$ws.Range("F11").Value = 2.148492139862341
$ws.Range("I11").Value = 1.703136487916495
$ws.Range("K11").Value = 0.9556792736053978
$ws.Range("L11").Value = 0.2399319355212981
$ws.Range("M11").Value = 0.2733192051428475
$ws.Range("N11").Value = 2.976232028045274
$ws.Range("B12").Value = 1.130262590177949
$ws.Range("C12").Value = 0.08700039087523237
$ws.Range("D12").Value = 0.115167299728494
$ws.Range("E12").Value = 0.05959108145640002
$ws.Range("F12").Value = 2.154056945733203
$ws.Range("I12").Value = 1.704636515718377
$ws.Range("K12").Value = 0.9696495702090715
$ws.Range("L12").Value = 0.2410024752727935
$ws.Range("M12").Value = 0.2756932812352204
$ws.Range("N12").Value = 2.973014223188926
$ws.Range("B13").Value = 1.127557202004141
$ws.Range("C13").Value = 0.08654071450109768
$ws.Range("D13").Value = 0.115174048320501
$ws.Range("E13").Value = 0.05959482748649236
$ws.Range("F13").Value = 2.152852006616072
$ws.Range("I13").Value = 1.704309041201483
$ws.Range("K13").Value = 0.9666378909718674
$ws.Range("L13").Value = 0.240771184058417
$ws.Range("M13").Value = 0.275181099190732
$ws.Range("N13").Value = 2.973702091713534
$ws.Range("B14").Value = 1.118747865356909
$ws.Range("C14").Value = 0.08504159471027606
$ws.Range("D14").Value = 0.1151960380203825
$ws.Range("E14").Value = 0.05960739905455981
$ws.Range("F14").Value = 2.148947049590106
$ws.Range("I14").Value = 1.703257905911357
$ws.Range("K14").Value = 0.9568272997885572
$ws.Range("L14").Value = 0.2400196798817404
$ws.Range("M14").Value = 0.2735141241756693
$ws.Range("N14").Value = 2.975964979803109
$ws.Range("B15").Value = 1.113362286660617
$ws.Range("C15").Value = 0.08412332734548045
$ws.Range("D15").Value = 0.1152094928617569
$ws.Range("E15").Value = 0.0596153723990458
$ws.Range("F15").Value = 2.146574060187774
$ws.Range("I15").Value = 1.702626987181148
$ws.Range("K15").Value = 0.9508265974651806
$ws.Range("L15").Value = 0.2395615037452359
$ws.Range("M15").Value = 0.2724956380074843
$ws.Range("N15").Value = 2.977366122231416
$ws.Range("B16").Value = 1.082659181756952
$ws.Range("C16").Value = 0.07886084363900636
$ws.Range("D16").Value = 0.115286377488502
$ws.Range("E16").Value = 0.05966525319135307
$ws.Range("F16").Value = 2.13326515652048
$ws.Range("I16").Value = 1.699209630462676
$ws.Range("K16").Value = 0.9165703510646779
$ws.Range("L16").Value = 0.2369687306962618
$ws.Range("M16").Value = 0.2666987458258845
$ws.Range("N16").Value = 2.985626207658072
$ws.Range("B17").Value = 1.063966284817809
$ws.Range("C17").Value = 0.07563217656398535
$ws.Range("D17").Value = 0.1153333503884895
$ws.Range("E17").Value = 0.05969959773082234
$ws.Range("F17").Value = 2.125359015237976
$ws.Range("I17").Value = 1.697290004949039
$ws.Range("K17").Value = 0.8956724927276696
$ws.Range("L17").Value = 0.235407482276031
$ws.Range("M17").Value = 0.263177947090611
$ws.Range("N17").Value = 2.990899261574242
$ws.Range("B18").Value = 1.053267045880858
$ws.Range("C18").Value = 0.07377499425244594
$ws.Range("D18").Value = 0.1153602979775155
$ws.Range("E18").Value = 0.05972073042496717
$ws.Range("F18").Value = 2.120906548548319
$ws.Range("I18").Value = 1.696250853091236
$ws.Range("K18").Value = 0.8836956949030537
$ws.Range("L18").Value = 0.2345202823214407
$ws.Range("M18").Value = 0.2611659026584192
$ws.Range("N18").Value = 2.994007758045925
$ws.Range("B19").Value = 1.049653477071757
$ws.Range("C19").Value = 0.07314615792333257
$ws.Range("D19").Value = 0.1153694100518585
$ws.Range("E19").Value = 0.05972812261995752
$ws.Range("F19").Value = 2.119415321625397
$ws.Range("I19").Value = 1.695910171001046
$ws.Range("K19").Value = 0.8796479602120542
$ws.Range("L19").Value = 0.2342217454880569
$ws.Range("M19").Value = 0.2604868975245864
$ws.Range("N19").Value = 2.995073217901734
$ws.Range("B20").Value = 1.065950751730924
$ws.Range("C20").Value = 0.07597588740762262
$ws.Range("D20").Value = 0.1153283573103252
$ws.Range("E20").Value = 0.05969579907263167
$ws.Range("F20").Value = 2.126190811001223
$ws.Range("I20").Value = 1.697487628358331
$ws.Range("K20").Value = 0.8978926452757605
$ws.Range("L20").Value = 0.2355725634291872
$ws.Range("M20").Value = 0.2635513943133603
$ws.Range("N20").Value = 2.990330113303656
$ws.Range("B21").Value = 1.121333215258403
$ws.Range("C21").Value = 0.08548192513575259
$ws.Range("D21").Value = 0.1151895821367006
$ws.Range("E21").Value = 0.05960364962342091
$ws.Range("F21").Value = 2.150090089173617
$ws.Range("I21").Value = 1.703563954741128
$ws.Range("K21").Value = 0.9597071227229605
$ws.Range("L21").Value = 0.2402399685594787
$ws.Range("M21").Value = 0.2740032166854078
$ws.Range("N21").Value = 2.975297176384458
$ws.Range("B22").Value = 1.15799426125011
$ws.Range("C22").Value = 0.09169400712855236
$ws.Range("D22").Value = 0.1150982390495301
$ws.Range("E22").Value = 0.05955564581912576
$ws.Range("F22").Value = 2.166555888893527
$ws.Range("I22").Value = 1.708113957990662
$ws.Range("K22").Value = 1.000489975969458
$ws.Range("L22").Value = 0.2433862761354106
$ws.Range("M22").Value = 0.2809497852379437
$ws.Range("N22").Value = 2.966146140201388
$ws.Range("B23").Value = 1.138385053603201
$ws.Range("C23").Value = 0.08837853512977745
$ws.Range("D23").Value = 0.1151470506891705
$ws.Range("E23").Value = 0.05958015077766277
$ws.Range("F23").Value = 2.157690309387405
$ws.Range("I23").Value = 1.705632564911454
$ws.Range("K23").Value = 0.9786883197214138
$ws.Range("L23").Value = 0.2416982676358401
$ws.Range("M23").Value = 0.2772316986358803
$ws.Range("N23").Value = 2.970968519723385
$ws.Range("B24").Value = 1.065053426090344
$ws.Range("C24").Value = 0.07582049879198394
$ws.Range("D24").Value = 0.1153306148588982
$ws.Range("E24").Value = 0.05969751212358876
$ws.Range("F24").Value = 2.125814466821296
$ws.Range("I24").Value = 1.697398081998067
$ws.Range("K24").Value = 0.89688879700347
$ws.Range("L24").Value = 0.2354978978981563
$ws.Range("M24").Value = 0.2633825211006879
$ws.Range("N24").Value = 2.990587185535958
$ws.Range("B25").Value = 0.9879580745888745
$ws.Range("C25").Value = 0.06226894508658631
$ws.Range("D25").Value = 0.1155259586164856
$ws.Range("E25").Value = 0.05987684571732199
$ws.Range("F25").Value = 2.095060163703309
$ws.Range("I25").Value = 1.690999427525355
$ws.Range("K25").Value = 0.810302290148968
$ws.Range("L25").Value = 0.2312812078828586
$ws.Range("M25").Value = 0.253739939928586
$ws.Range("N25").Value = 3.006150687704675
